$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.07002266666666
$ws.Range("H2").Value = 36.21006799999999
$ws.Range("I2").Value = 0.7601982364861632
$ws.Range("J2").Value = 0.7601982364861634
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.1030763333333333
$ws.Range("N2").Value = 0.309229
$ws.Range("O2").Value = 0.01126512502660735
$ws.Range("P2").Value = 0.01126512502660735
$ws.Range("Q2").Value = 1.244133679730222
$ws.Range("R2").Value = 11.197203117572
$ws.Range("S2").Value = 0.008563728179023053
$ws.Range("T2").Value = 0.008563728179023055

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.07002266666666
$ws.Range("H3").Value = 36.21006799999999
$ws.Range("I3").Value = 0.7601982364861632
$ws.Range("J3").Value = 0.7601982364861634
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.046962666666667
$ws.Range("N3").Value = 27.140888
$ws.Range("O3").Value = 0.9887348749733926
$ws.Range("P3").Value = 0.9887348749733927
$ws.Range("Q3").Value = 109.1970444511538
$ws.Range("R3").Value = 982.7734000603838
$ws.Range("S3").Value = 0.7516345083071402
$ws.Range("T3").Value = 0.7516345083071404

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.308268
$ws.Range("H4").Value = 3.924804
$ws.Range("I4").Value = 0.08239777620284613
$ws.Range("J4").Value = 0.08239777620284613
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1030763333333333
$ws.Range("N4").Value = 0.309229
$ws.Range("O4").Value = 0.01126512502660735
$ws.Range("P4").Value = 0.01126512502660735
$ws.Range("Q4").Value = 0.1348514684573333
$ws.Range("R4").Value = 1.213663216116
$ws.Range("S4").Value = 0.0009282212508394738
$ws.Range("T4").Value = 0.0009282212508394739

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.308268
$ws.Range("H5").Value = 3.924804
$ws.Range("I5").Value = 0.08239777620284613
$ws.Range("J5").Value = 0.08239777620284613
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.046962666666667
$ws.Range("N5").Value = 27.140888
$ws.Range("O5").Value = 0.9887348749733926
$ws.Range("P5").Value = 0.9887348749733927
$ws.Range("Q5").Value = 11.83585175399467
$ws.Range("R5").Value = 106.522665785952
$ws.Range("S5").Value = 0.08146955495200665
$ws.Range("T5").Value = 0.08146955495200667

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.499176666666667
$ws.Range("H6").Value = 7.49753
$ws.Range("I6").Value = 0.1574039873109905
$ws.Range("J6").Value = 0.1574039873109906
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.1030763333333333
$ws.Range("N6").Value = 0.309229
$ws.Range("O6").Value = 0.01126512502660735
$ws.Range("P6").Value = 0.01126512502660735
$ws.Range("Q6").Value = 0.2576059671522222
$ws.Range("R6").Value = 2.31845370437
$ws.Range("S6").Value = 0.001773175596744826
$ws.Range("T6").Value = 0.001773175596744826

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.499176666666667
$ws.Range("H7").Value = 7.49753
$ws.Range("I7").Value = 0.1574039873109905
$ws.Range("J7").Value = 0.1574039873109906
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.046962666666667
$ws.Range("N7").Value = 27.140888
$ws.Range("O7").Value = 0.9887348749733926
$ws.Range("P7").Value = 0.9887348749733927
$ws.Range("Q7").Value = 22.60995800073778
$ws.Range("R7").Value = 203.48962200664
$ws.Range("S7").Value = 0.1556308117142457
$ws.Range("T7").Value = 0.1556308117142458
